$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "41+1=42"
$t.Cell(1, 2).Range.Text = "15+22=37"
$t.Cell(1, 3).Range.Text = "60-53=7"
$t.Cell(1, 4).Range.Text = "13+1=14"
$t.Cell(1, 5).Range.Text = "16+36=52"
$t.Cell(2, 1).Range.Text = "87-36=51"
$t.Cell(2, 2).Range.Text = "36-13=23"
$t.Cell(2, 3).Range.Text = "3+32=35"
$t.Cell(2, 4).Range.Text = "8-1=7"
$t.Cell(2, 5).Range.Text = "80-57=23"
$t.Cell(3, 1).Range.Text = "7+46=53"
$t.Cell(3, 2).Range.Text = "99-78=21"
$t.Cell(3, 3).Range.Text = "71-64=7"
$t.Cell(3, 4).Range.Text = "5+45=50"
$t.Cell(3, 5).Range.Text = "32-15=17"
$t.Cell(4, 1).Range.Text = "46+51=97"
$t.Cell(4, 2).Range.Text = "88-27=61"
$t.Cell(4, 3).Range.Text = "67+29=96"
$t.Cell(4, 4).Range.Text = "99-27=72"
$t.Cell(4, 5).Range.Text = "44-29=15"
$t.Cell(5, 1).Range.Text = "13+29=42"
$t.Cell(5, 2).Range.Text = "19-7=12"
$t.Cell(5, 3).Range.Text = "82-3=79"
$t.Cell(5, 4).Range.Text = "53+39=92"
$t.Cell(5, 5).Range.Text = "89-60=29"
$t.Cell(6, 1).Range.Text = "55-45=10"
$t.Cell(6, 2).Range.Text = "85-82=3"
$t.Cell(6, 3).Range.Text = "22+47=69"
$t.Cell(6, 4).Range.Text = "24-4=20"
$t.Cell(6, 5).Range.Text = "27-6=21"
$t.Cell(7, 1).Range.Text = "3+91=94"
$t.Cell(7, 2).Range.Text = "44+40=84"
$t.Cell(7, 3).Range.Text = "10-9=1"
$t.Cell(7, 4).Range.Text = "27+34=61"
$t.Cell(7, 5).Range.Text = "19-17=2"
$t.Cell(8, 1).Range.Text = "0+29=29"
$t.Cell(8, 2).Range.Text = "79+10=89"
$t.Cell(8, 3).Range.Text = "28+46=74"
$t.Cell(8, 4).Range.Text = "50+33=83"
$t.Cell(8, 5).Range.Text = "62-1=61"
$t.Cell(9, 1).Range.Text = "29+18=47"
$t.Cell(9, 2).Range.Text = "82-68=14"
$t.Cell(9, 3).Range.Text = "25+70=95"
$t.Cell(9, 4).Range.Text = "34+41=75"
$t.Cell(9, 5).Range.Text = "56-17=39"
$t.Cell(10, 1).Range.Text = "4+33=37"
$t.Cell(10, 2).Range.Text = "95-26=69"
$t.Cell(10, 3).Range.Text = "39-10=29"
$t.Cell(10, 4).Range.Text = "4+63=67"
$t.Cell(10, 5).Range.Text = "21-14=7"
$t.Cell(11, 1).Range.Text = "85-32=53"
$t.Cell(11, 2).Range.Text = "75-6=69"
$t.Cell(11, 3).Range.Text = "39+2=41"
$t.Cell(11, 4).Range.Text = "21+18=39"
$t.Cell(11, 5).Range.Text = "69-19=50"
$t.Cell(12, 1).Range.Text = "50-42=8"
$t.Cell(12, 2).Range.Text = "47-21=26"
$t.Cell(12, 3).Range.Text = "72+21=93"
$t.Cell(12, 4).Range.Text = "89+8=97"
$t.Cell(12, 5).Range.Text = "69+28=97"
$t.Cell(13, 1).Range.Text = "62-7=55"
$t.Cell(13, 2).Range.Text = "88-34=54"
$t.Cell(13, 3).Range.Text = "72-27=45"
$t.Cell(13, 4).Range.Text = "89+3=92"
$t.Cell(13, 5).Range.Text = "46-38=8"
$t.Cell(14, 1).Range.Text = "7+47=54"
$t.Cell(14, 2).Range.Text = "98-37=61"
$t.Cell(14, 3).Range.Text = "76+14=90"
$t.Cell(14, 4).Range.Text = "31+9=40"
$t.Cell(14, 5).Range.Text = "69-12=57"
$t.Cell(15, 1).Range.Text = "62-12=50"
$t.Cell(15, 2).Range.Text = "72-44=28"
$t.Cell(15, 3).Range.Text = "22+25=47"
$t.Cell(15, 4).Range.Text = "99-64=35"
$t.Cell(15, 5).Range.Text = "10+63=73"
$t.Cell(16, 1).Range.Text = "91-78=13"
$t.Cell(16, 2).Range.Text = "38+12=50"
$t.Cell(16, 3).Range.Text = "52-42=10"
$t.Cell(16, 4).Range.Text = "80-11=69"
$t.Cell(16, 5).Range.Text = "56-21=35"
$t.Cell(17, 1).Range.Text = "91-90=1"
$t.Cell(17, 2).Range.Text = "60+4=64"
$t.Cell(17, 3).Range.Text = "31-3=28"
$t.Cell(17, 4).Range.Text = "54-20=34"
$t.Cell(17, 5).Range.Text = "36+41=77"
$t.Cell(18, 1).Range.Text = "58+9=67"
$t.Cell(18, 2).Range.Text = "34-9=25"
$t.Cell(18, 3).Range.Text = "89-41=48"
$t.Cell(18, 4).Range.Text = "69-68=1"
$t.Cell(18, 5).Range.Text = "76-73=3"
$t.Cell(19, 1).Range.Text = "6-5=1"
$t.Cell(19, 2).Range.Text = "58+5=63"
$t.Cell(19, 3).Range.Text = "51+8=59"
$t.Cell(19, 4).Range.Text = "36-6=30"
$t.Cell(19, 5).Range.Text = "25+41=66"
$t.Cell(20, 1).Range.Text = "55+34=89"
$t.Cell(20, 2).Range.Text = "14-3=11"
$t.Cell(20, 3).Range.Text = "12+87=99"
$t.Cell(20, 4).Range.Text = "97+2=99"
$t.Cell(20, 5).Range.Text = "61-6=55"
